# "clean up blocks 1-4"
#
# Slide 5 ("Getting Started"), the subtitle placeholder shape, has a
# paragraph that ends with a sentence pointing people at the session
# webpage, followed by a separate yellow-highlighted "insert link"
# placeholder run, and then an otherwise-empty paragraph used only for
# extra line spacing. Clean this up:
#   - reword the sentence so it reads "...which can be found on the
#     session webpage." instead of "...session webpage: [insert link]"
#   - drop the yellow highlight placeholder run entirely
#   - fold the now-pointless empty paragraph into this one (it only
#     existed to carry the trailing endParaRPr / extra spacing)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

$para = $tr.Paragraphs(3, 1)

# 1) Trim the first run back to "...found on the " (drop the trailing
#    "session webpage: " that used to lead into the highlighted run).
$run1 = $para.Runs(1, 1)
$run1.Text = "For this course there is a template set up for you to fork (copy), which can be found on the "

# 2) Remove the whole highlighted "insert link" run and type a plain
#    "session webpage." run in its place.
$para = $tr.Paragraphs(3, 1)
$run2 = $para.Runs(2, 1)
$hlStart = $run2.Start
$hlLen = $run2.Text.Length
$tr.Characters($hlStart, $hlLen).Delete()

$para = $tr.Paragraphs(3, 1)
$run1 = $para.Runs(1, 1)
[void]$run1.InsertAfter("session webpage.")

# 3) The next paragraph is now redundant (it only held extra 150%
#    line-spacing with no visible text) - delete it so its end-of-paragraph
#    run properties collapse into the paragraph above.
$para = $tr.Paragraphs(3, 1)
$nextPara = $tr.Paragraphs(4, 1)
[void]$nextPara.Delete()
